$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 4 label cell so the new row 5 label (A5) matches
# the existing year-label style (bold, bordered, centered) used by A2:A4.
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A5").Value = "2021年"
$ws.Range("B5").Value = 53.545
$ws.Range("C5").Value = 65.36199999999999
$ws.Range("D5").Value = 35.767
$ws.Range("E5").Value = 39.173
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = 24.103
$ws.Range("H5").Value = 37.952
$ws.Range("I5").Value = 41.557
$ws.Range("J5").Value = 41.987
$ws.Range("K5").Value = 47.434
$ws.Range("L5").Value = 30.615
$ws.Range("M5").Value = 32.956
$ws.Range("N5").Value = 18.922
$ws.Range("O5").Value = 8.491
$ws.Range("P5").Value = 34.388
$ws.Range("Q5").Value = 32.787
$ws.Range("R5").Value = 9.207000000000001
$ws.Range("S5").Value = 17.021
$ws.Range("T5").Value = 37.871
$ws.Range("U5").Value = 4.527
$ws.Range("V5").Value = 49.232
$ws.Range("W5").Value = 50.794
$ws.Range("X5").Value = 2.235
$ws.Range("Y5").Value = 3.264
$ws.Range("Z5").Value = 3.394
$ws.Range("AA5").Value = 3.108
$ws.Range("AB5").Value = 52.274
$ws.Range("AC5").Value = 29.91
$ws.Range("AD5").Value = 21.272
$ws.Range("AE5").Value = 4.196
$ws.Range("AF5").Value = 31.327
$ws.Range("AG5").Value = 20.388
$ws.Range("AH5").Value = 58.245
$ws.Range("AI5").Value = 50.907
$ws.Range("AJ5").Value = 28.597
$ws.Range("AK5").Value = 29.721
$ws.Range("AL5").Value = 5.683
$ws.Range("AM5").Value = 20.376
$ws.Range("AN5").Value = 33.843
$ws.Range("AO5").Value = 46.195
$ws.Range("AP5").Value = 21.871
$ws.Range("AQ5").Value = 9.361000000000001
$ws.Range("AR5").Value = 38.043
$ws.Range("AS5").Value = 24.38
$ws.Range("AT5").Value = 3.867
